$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, taken from the updated cryptos feed.
# A handful of Price-column cells are purely numeric-looking
# ("134.80", "0.0220", ...); those are prefixed with a literal
# apostrophe so Excel keeps them as text (matching the source
# data, which stores every cell as a string) instead of parsing
# them into numbers and silently dropping significant trailing
# zeros (134.80 -> 134.8, 0.0220 -> 0.022, etc).
$updates = [ordered]@{
    'D2' = '60.050.51'
    'E2' = '  +0.23%  '
    'D3' = '2.399.68'
    'E3' = '  -0.72%  '
    'E4' = '  +0.05%  '
    'D5' = "'559.67"
    'E5' = '  +1.44%  '
    'D6' = "'134.80"
    'E6' = '  -1.97%  '
    'E7' = '  +0.07%  '
    'D8' = "'0.587"
    'E8' = '  -0.25%  '
    'E9' = '  +0.16%  '
    'E10' = '  -0.59%  '
    'E11' = '  +1.34%  '
    'D12' = "'0.345"
    'E12' = '  -2.71%  '
    'D13' = "'24.60"
    'E13' = '  -3.59%  '
    'D14' = '2.828.30'
    'E14' = '  -0.67%  '
    'D15' = '60.016.15'
    'E15' = '  +0.30%  '
    'D16' = "'0.0000137"
    'E16' = '  -0.16%  '
    'D17' = '2.401.17'
    'E17' = '  -0.85%  '
    'D18' = "'11.13"
    'E18' = '  -2.04%  '
    'D19' = "'4.52"
    'E19' = '  +2.45%  '
    'D20' = "'322.85"
    'E20' = '  -1.98%  '
    'D21' = "'6.77"
    'E21' = '  +1.64%  '
    'E22' = '  +0.06%  '
    'D23' = "'64.17"
    'E23' = '  -3.46%  '
    'D24' = "'0.173"
    'E24' = '  +0.02%  '
    'D25' = "'8.48"
    'E25' = '  -2.16%  '
    'E26' = '  -0.01%  '
    'D27' = "'1.39"
    'E27' = '  +0.72%  '
    'E28' = '  +1.82%  '
    'D29' = '0.0₃0767'
    'E29' = '  -1.38%  '
    'D30' = "'170.99"
    'E30' = '  +1.05%  '
    'D31' = "'6.11"
    'E31' = '  -0.16%  '
    'E32' = '  +6.72%  '
    'D33' = "'0.401"
    'E33' = '  -1.98%  '
    'D34' = "'18.25"
    'E34' = '  -2.42%  '
    'E35' = '  +0.04%  '
    'E36' = '  +2.21%  '
    'E37' = '  +0.04%  '
    'D38' = "'4.17"
    'E38' = '  -1.37%  '
    'D39' = "'323.27"
    'E39' = '  +2.81%  '
    'D40' = "'1.59"
    'E40' = '  -0.77%  '
    'D41' = "'38.69"
    'E41' = '  -2.23%  '
    'D42' = "'147.35"
    'E42' = '  +6.34%  '
    'D43' = "'3.55"
    'D44' = "'0.0967"
    'E44' = '  +0.09%  '
    'D45' = "'19.86"
    'E45' = '  +1.36%  '
    'D46' = "'0.0514"
    'E46' = '  -1.15%  '
    'D47' = "'0.575"
    'E47' = '  -0.72%  '
    'D48' = "'0.0220"
    'E48' = '  -1.89%  '
    'D50' = "'1.56"
    'E50' = '  -1.12%  '
    'E51' = '  +0.22%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Host "Updated $($updates.Count) cells"
